# Edit script for LOT2038.docx
$d = $word.ActiveDocument
$wdReplaceNone = 0
$wdFindContinue = 1
$wdReplaceOne = 1

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceOne)
    if (-not $ok) {
        Write-Host "NOT FOUND: $old"
    }
}

# 1. Ativacao date
Replace-Text "Ativação: 01/01/2018" "Ativação: 01/01/2025"

# 2. Fill empty italic "Objetivos" English run (empty paragraph right after the
#    Portuguese objectives paragraph ending in "... destilação, maturação.")
$pars = $d.Paragraphs
for ($i = 1; $i -le $pars.Count; $i++) {
    $p = $pars.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "" -and $p.Range.Font.Italic -and $i -gt 1) {
        $prevTxt = $pars.Item($i - 1).Range.Text
        if ($prevTxt -like "*destilação, maturação.*") {
            $p.Range.InsertAfter("Provide the student with theoretical knowledge in the technological processes of preparing fermented and distilled beverages. Knowledge of raw materials, wort preparation, fermentation technology, distillation, beverage maturation.")
            break
        }
    }
}

# 3. Programa resumido - Portuguese
Replace-Text "Generalidades, tipos de bebidas e preparação do mosto; tipos de leveduras; fermentação alcoólica. Acabamento: maturação, filtração, destilação, envelhecimento. Análise química e sensorial." "Generalidades, tipos de bebidas alcoólicas e preparação do mosto; tipos de leveduras; fermentação alcoólica. Acabamento: maturação, filtração, destilação, envelhecimento. Bebidas fermentadas não-alcoólicas: generalidades e processo produtivo. Análise química e sensorial."

# 4. Programa resumido - English (italic)
Replace-Text "Generalities, types of beverages and wort preparation; types of yeasts; alcoholic fermentation. Finishing: maturation, filtration, distillation, aging. Chemical and sensorial analyzes." "General information, types of alcoholic beverages and must preparation; types of yeast; alcoholic fermentation. Finishing: maturation, filtration, distillation, aging. Non-alcoholic fermented beverages: generalities and production process. Chemical and sensory analysis of beverages."

# 5. Programa - Portuguese
Replace-Text "1. Generalidades: origem das bebidas; matériasprimas.2. Tipos de bebidas: bebidas fermentadas (cervejas, fermentados de frutas, vinhos, sidras),bebidas destiladas (aguardentes, destilados de vinhos, graspa, pisco, rum, tequila, tiquira,uísque), bebidas retificadas (vodka, gim), bebidas obtidas por misturas (licores, sangria, cooler).3. Preparação do mosto: pé de cuba, características físicoquímicas,correção do mosto.4. Tipos de leveduras: leveduras selvagens, leveduras mistas, leveduras selecionadas.5. Fermentação alcoólica: controle da fermentação, rendimento da fermentação, produtossecundários.6. Acabamento: controle da maturação, destilação em alambiques e em colunas, determinação dograu alcoólico, armazenamento, tipos de madeiras, cor, volume, composição da bebida, legislação.7. Análise química: composição da bebida, legislação.8. Análise sensorial: aromas das bebidas e aceitação." "1. Generalidades: origem das bebidas; matérias-primas. 2. Tipos de bebidas: bebidas fermentadas (cervejas, fermentados de frutas, vinhos, sidras), bebidas destiladas (aguardentes, destilados de vinhos, graspa, pisco, rum, tequila, tiquira, uísque), bebidas retificadas (vodka, gim), bebidas obtidas por misturas (licores, sangria, cooler). 3. Preparação do mosto: pé de cuba, características físicoquímicas, correção do mosto. 4. Tipos de leveduras: leveduras selvagens, leveduras mistas, leveduras selecionadas. 5. Fermentação alcoólica: controle da fermentação, rendimento da fermentação, produtos secundários. 6. Acabamento: controle da maturação, destilação em alambiques e em colunas, determinação do grau alcoólico, armazenamento, tipos de madeiras, cor, volume, composição da bebida, legislação. 7. Preparação de bebiidas fermentadas não alcoólicas (iogurte, leite fermentado, kefir e outras). 8.Análise química: composição da bebida, legislação. 8. Análise sensorial: aromas das bebidas e aceitação."

# 6. Programa - English (italic)
Replace-Text "1. Generalities: beverages origin; rawmaterial.2. Types of beverages: fermented beverages (beer, fruits fermented fruits, wines, ciders),distilled beverages (spirits, wines distilled, graspa, pisco, rum, tequila, whisky), rectifiedbeverages (vodka, gin), beverages obtained from mixtures (liquors, sangria, cooler).3. Worth preparation: preparation of the inoculum, physicalchemicalcharacteristics, correction ofthe wort.4. Types of yeasts: wild yeasts, mixed yeasts, selected yeasts.5. Alcoholic fermentation: fermentation control, fermentation yield, secondary products.6. Finishing: maturation control, distillation in stills and in columns, alcoholic degreedetermination, storage, types of woods, color, volume, beverage composition, legislation.7. Chemical analysis: beverage composition, legislation.8. Sensorial analysis: beverages flavor and acceptance" "1. General: origin of the beverages; 2. Types of beverages: fermented beverages (beers, fermented fruits, wines, ciders), distilled beverages (cachaça, wine distillates, graspa, pisco, rum, tequila, tiquira, whiskey), rectified beverages (vodka, gin), beverages obtained by mixing (liqueurs, sangria, cooler). 3. Preparation of the must: foot of vat, physicochemical characteristics, correction of the wort. 4. Types of yeast: wild yeast, mixed yeast, selected yeast. 5. Alcoholic fermentation: fermentation control, fermentation yield, secondary products. 6. Finishing: control of maturation, distillation in stills and columns, determination of alcoholic content, storage, types of wood, color, volume, composition of the drink, legislation. 7. Preparation of non-alcoholic fermented beverages (yogurt, fermented milk, kefir and others). 8. Chemical analysis: composition of the beverages, legislation. 9. Sensory analysis: beverage aromas and acceptance"

# 7. Bibliografia
Replace-Text "1. AQUARONE, E.; BORZANI, W.; SCHMIDELL, W.; LIMA, U. A. Biotecnologia na Produção deAlimentos. V. 4, Biotecnologia Industrial, São Paulo: Edgard Blücher Ltda. 2001.3. DUVAL, G. Fabricação de Vinhos de Frutas. S.I.A.RJ:Ministério da Agricultura, 1947.4. HOUGH, J.S. Biotecnología de La cerveza y de la malta. Editorial ACRIBA S/A, 1978.5. LIMA, U. A. Aguardente: fabricação em pequenas destilarias. Ed. FEALQ. 1999.6. MARTINELLI FILHO, A. Tecnologia de Vinhos e Vinagres de Frutas. Agroindústria de BaixoInvestimento. Departamento de Tecnologia Rural da ESALQ/USP.7. MORRETO, E. et al. Vinhos e Vinagres: Processamento e Análises. FlorianópolisEditoraUFSC, 1988.8. PACHECO, A. O. Manual do Bar. São Paulo. Editora SENAC, 1996.9. STANIER, R. Y.; INGRAHAM, J. L., WHEELIS, M. L.; PAINTER, P. R. The Microbial World.Englewood Cliffs, New Jersey, 1986.10.Venturini Filho, W.G. Bebidas Alcoólicas. Ciência e Tecnologia. São Paulo. Edgar Blucher Ltda. 2a. Edição. 2016. 575 p." "1) DA SILVA, N., JUNQUEIRA, V. C. A., DE ARRUDA SILVEIRA, N. F., TANIWAKI, M. H., GOMES, R. A. R., OKAZAKI, M. M. Manual de métodos de análise microbiológica de alimentos e água. Editora Blucher, 2017. 2) DA-SILVA, R.; LAGO-VANZELA, E. S.; BAFFI, M. A. Uvas e vinhos: química, bioquímica e microbiologia. São Paulo, Editora Senac, 2015. 3) DE OLIVEIRA MORAES, I. Biotecnologia Industrial: biotecnologia na produção de alimentos. Vol. 4. 2ª Ed. Editora Blucher, 2021. 4) MARTIN, J. G. P., DE DEA LINDNER, J. Microbiologia de alimentos fermentados. Editora Blucher, 2022. 5) MENEZES e SILVA, C.H.P. Microbiologia da cerveja - Do básico ao avançado, o guia definitivo. Editora LF, 2019. 6) MUXEL, A. A. Química da Cerveja: Uma Abordagem Química e Bioquímica das Matérias-Primas, Processo de Produção e da Composição dos Compostos de Sabores da Cerveja. Editora Appris, 2022. 7) VENTURINI FILHO, W. G. Bebidas alcoólicas: ciência e tecnologia. Vol. 1. Editora Blucher, 2021."
